$d = $word.ActiveDocument

# Rename the Pearson logo images in the footers: image2.png -> image1.png
$d.Sections.Item(1).Footers.Item(1).Range.InlineShapes.Item(1).Name = "image1.png"
$d.Sections.Item(1).Footers.Item(2).Range.InlineShapes.Item(1).Name = "image1.png"

# Rename the BTec logo images in the headers: image1.jpg -> image2.jpg
$d.Sections.Item(1).Headers.Item(1).Range.InlineShapes.Item(1).Name = "image2.jpg"
$d.Sections.Item(1).Headers.Item(2).Range.InlineShapes.Item(1).Name = "image2.jpg"
